$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.332.47'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.796.01'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.73'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.33'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.70%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.794.26'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.516'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.48%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.98'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +10.74%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.04'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.431.52'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.824.46'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.326.36'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.07'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.36%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.62%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '461.81'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.57'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.15%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.58%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.49'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.98'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.10'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.98'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.947.39'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -5.83%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.23'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.01'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.02'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.96%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.148'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +7.22%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.976'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.18'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '46.99'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.54%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '152.99'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.296'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '42.88'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.40'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.36'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.85'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.13'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.57%  '
